# Apply the "NewSell" template edit:
#   - A2 "Lapiz Grafito" -> "Cuaderno", B2 5 -> 20 (C2 price stays 100)
#   - New row 3: "Libreta " (trailing space) / 13 / 55
#   - Active cell / selection moves to D6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing product row (Cuaderno replaces Lapiz Grafito, quantity changes)
$ws.Range("A2").Value = "Cuaderno"
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 100

# Add the new product row
$ws.Range("A3").Value = "Libreta "
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 55

# Leave the selection where the user last clicked
$null = $ws.Range("D6").Select()
